# Generate Report for Handoff
#
# Updates the localization-status workbook so that the overview/status rows
# reflect that the content is now "Ready for handoff" (instead of "Handed
# back: in sync with en-US"), refreshes the associated timestamps, and
# shrinks the now-shorter "Status" columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Refresh timestamps for the new handoff
# Overview "Latest HO Xliff Generate Date" column (G) for both locales shared
# the same timestamp; de-de's handoff datetime mirrored it too.
$wsOverview.Range("G2").Value = "2016-09-05 23:13:17"
$wsDeDe.Range("H2").Value = "2016-09-05 23:13:17"

# zh-cn's "Latest Handoff Datetime" moved to a slightly earlier moment.
$wsZhCn.Range("H2").Value = "2016-09-05 23:13:13"

# --- Narrow the Status columns now that the text is shorter.
$newStatusWidth = 98 / 6   # closest attainable ColumnWidth to the target fit width

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth
